$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.471.58"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "3.974.88"
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.95"
$ws.Range("E5").Value = "  +2.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.64"
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("D7").Value = "3.967.99"
$ws.Range("E7").Value = "  -2.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.688"
$ws.Range("E8").Value = "  -3.93%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("E10").Value = "  -3.56%  "

$ws.Range("E11").Value = "  -7.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.40"
$ws.Range("E12").Value = "  +13.51%  "

$ws.Range("E13").Value = "  -6.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.73"
$ws.Range("E14").Value = "  -3.04%  "

$ws.Range("D15").Value = "4.610.21"
$ws.Range("E15").Value = "  -2.44%  "

$ws.Range("D16").Value = "3.978.30"
$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  -3.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.65"
$ws.Range("E18").Value = "  -3.03%  "

$ws.Range("E19").Value = "  -1.72%  "

$ws.Range("E20").Value = "  -5.31%  "

$ws.Range("D21").Value = "71.331.90"
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "433.30"
$ws.Range("E22").Value = "  -3.81%  "

$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "97.38"
$ws.Range("E24").Value = "  -6.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.25"
$ws.Range("E25").Value = "  +4.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.62"
$ws.Range("E26").Value = "  -2.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.13"
$ws.Range("E27").Value = "  +24.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.40"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("E29").Value = "  -3.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.92"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.83"
$ws.Range("E31").Value = "  -3.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  +18.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.49"
$ws.Range("E33").Value = "  +20.71%  "

$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.41"
$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "677.33"
$ws.Range("E36").Value = "  -1.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.81"
$ws.Range("E37").Value = "  -2.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.446"
$ws.Range("E38").Value = "  +2.38%  "

$ws.Range("D39").Value = "0.0₃0825"
$ws.Range("E39").Value = "  -9.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.149"
$ws.Range("E40").Value = "  -3.23%  "

$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("E44").Value = "  -3.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.24"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.42"
$ws.Range("E46").Value = "  +4.75%  "

$ws.Range("E47").Value = "  -4.92%  "

$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("E51").Value = "  -8.97%  "
